$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Mondaay"
$ws.Range("D2").Value = "mondaay"
$ws.Range("E2").Value = "MONDAAY"
$ws.Range("F2").Value = 7
